# Applies the scheduled-runner profit recalculations to the Durandal_Profits workbook.
# Each FFXIV crafting-class sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) gets updated
# currentAveragePrice / LevePrice / LeveProfit columns (H-N) for a handful of rows,
# reflecting refreshed market-board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 29.5625
$ws.Range("I11").Value = 29.5625
$ws.Range("K11").Value = 29.5625
$ws.Range("M11").Value = 110.4375

$ws.Range("H40").Value = 1419.579
$ws.Range("I40").Value = 1371.5714
$ws.Range("J40").Value = 1554
$ws.Range("K40").Value = 1371.5714
$ws.Range("L40").Value = 1554
$ws.Range("M40").Value = -1196.5714
$ws.Range("N40").Value = -1904

$ws.Range("H51").Value = 2538.6667
$ws.Range("J51").Value = 2701.2727
$ws.Range("L51").Value = 2701.2727
$ws.Range("N51").Value = -3669.2727

$ws.Range("H64").Value = 3783.3333
$ws.Range("I64").Value = 3400
$ws.Range("J64").Value = 4166.6665
$ws.Range("K64").Value = 3400
$ws.Range("L64").Value = 4166.6665
$ws.Range("M64").Value = -3152
$ws.Range("N64").Value = -4662.6665

$ws.Range("H67").Value = 3783.3333
$ws.Range("I67").Value = 3400
$ws.Range("J67").Value = 4166.6665
$ws.Range("K67").Value = 3400
$ws.Range("L67").Value = 4166.6665
$ws.Range("M67").Value = -2542
$ws.Range("N67").Value = -5882.6665

$ws.Range("H112").Value = 6499.316
$ws.Range("J112").Value = 6499.316
$ws.Range("L112").Value = 19497.948
$ws.Range("N112").Value = -21713.948

$ws.Range("H121").Value = 964.8261
$ws.Range("J121").Value = 922.05
$ws.Range("L121").Value = 2766.15
$ws.Range("N121").Value = -6260.15

$ws.Range("H135").Value = 2989.4546
$ws.Range("I135").Value = 3901
$ws.Range("J135").Value = 558.6667
$ws.Range("K135").Value = 35109
$ws.Range("L135").Value = 5028.0003
$ws.Range("M135").Value = -32574
$ws.Range("N135").Value = -10098.0003

$ws.Range("H138").Value = 3952.186
$ws.Range("I138").Value = 1831.6666
$ws.Range("J138").Value = 6630.737
$ws.Range("K138").Value = 5494.9998
$ws.Range("L138").Value = 19892.211
$ws.Range("M138").Value = -354.9997999999996
$ws.Range("N138").Value = -30172.211

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1956.683
$ws.Range("I61").Value = 1640.1936
$ws.Range("J61").Value = 2937.8
$ws.Range("K61").Value = 1640.1936
$ws.Range("L61").Value = 2937.8
$ws.Range("M61").Value = -1428.1936
$ws.Range("N61").Value = -3361.8

$ws.Range("H136").Value = 1956.683
$ws.Range("I136").Value = 1640.1936
$ws.Range("J136").Value = 2937.8
$ws.Range("K136").Value = 4920.5808
$ws.Range("L136").Value = 8813.400000000001
$ws.Range("M136").Value = -2370.5808
$ws.Range("N136").Value = -13913.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 127188.875
$ws.Range("I86").Value = 2301.2
$ws.Range("K86").Value = 2301.2
$ws.Range("M86").Value = -1178.2

$ws.Range("H89").Value = 127188.875
$ws.Range("I89").Value = 2301.2
$ws.Range("K89").Value = 11506
$ws.Range("M89").Value = -5890

$ws.Range("H99").Value = 2202.2173
$ws.Range("I99").Value = 1278.3334
$ws.Range("K99").Value = 1278.3334
$ws.Range("M99").Value = 219.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 817.11365
$ws.Range("I58").Value = 729.40625
$ws.Range("J58").Value = 1051
$ws.Range("K58").Value = 729.40625
$ws.Range("L58").Value = 1051
$ws.Range("M58").Value = -526.40625
$ws.Range("N58").Value = -1457

$ws.Range("H136").Value = 817.11365
$ws.Range("I136").Value = 729.40625
$ws.Range("J136").Value = 1051
$ws.Range("K136").Value = 2188.21875
$ws.Range("L136").Value = 3153
$ws.Range("M136").Value = 361.78125
$ws.Range("N136").Value = -8253

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 539.75
$ws.Range("I107").Value = 388.85294
$ws.Range("J107").Value = 906.2143
$ws.Range("K107").Value = 1166.55882
$ws.Range("L107").Value = 2718.6429
$ws.Range("M107").Value = 753.44118
$ws.Range("N107").Value = -6558.6429

$ws.Range("H131").Value = 5264201
$ws.Range("I131").Value = 1037.5714
$ws.Range("J131").Value = 5682862
$ws.Range("K131").Value = 3112.7142
$ws.Range("L131").Value = 17048586
$ws.Range("M131").Value = 1927.2858
$ws.Range("N131").Value = -17058666

$ws.Range("H137").Value = 10868.25
$ws.Range("I137").Value = 15502.714
$ws.Range("J137").Value = 4380
$ws.Range("K137").Value = 46508.142
$ws.Range("L137").Value = 13140
$ws.Range("M137").Value = -41408.142
$ws.Range("N137").Value = -23340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 28000
$ws.Range("J64").Value = 28000
$ws.Range("L64").Value = 28000
$ws.Range("N64").Value = -28496

$ws.Range("H67").Value = 28000
$ws.Range("J67").Value = 28000
$ws.Range("L67").Value = 28000
$ws.Range("N67").Value = -29716

$ws.Range("H74").Value = 59800
$ws.Range("J74").Value = 59800
$ws.Range("L74").Value = 59800
$ws.Range("N74").Value = -61672

$ws.Range("H77").Value = 59800
$ws.Range("J77").Value = 59800
$ws.Range("L77").Value = 179400
$ws.Range("N77").Value = -188760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1011.8
$ws.Range("I22").Value = 463.66666
$ws.Range("J22").Value = 1148.8334
$ws.Range("K22").Value = 463.66666
$ws.Range("L22").Value = 1148.8334
$ws.Range("M22").Value = -168.66666
$ws.Range("N22").Value = -1738.8334

$ws.Range("H27").Value = 1011.8
$ws.Range("I27").Value = 463.66666
$ws.Range("J27").Value = 1148.8334
$ws.Range("K27").Value = 463.66666
$ws.Range("L27").Value = 1148.8334
$ws.Range("M27").Value = -356.66666
$ws.Range("N27").Value = -1362.8334

$ws.Range("H62").Value = 12000
$ws.Range("J62").Value = 12000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13248

$ws.Range("H65").Value = 12000
$ws.Range("J65").Value = 12000
$ws.Range("L65").Value = 36000
$ws.Range("N65").Value = -42240

$ws.Range("H68").Value = 1851.2727
$ws.Range("I68").Value = 1595.7715
$ws.Range("J68").Value = 2844.889
$ws.Range("K68").Value = 1595.7715
$ws.Range("L68").Value = 2844.889
$ws.Range("M68").Value = -846.7715000000001
$ws.Range("N68").Value = -4342.889

$ws.Range("H71").Value = 1851.2727
$ws.Range("I71").Value = 1595.7715
$ws.Range("J71").Value = 2844.889
$ws.Range("K71").Value = 7978.8575
$ws.Range("L71").Value = 14224.445
$ws.Range("M71").Value = -4234.8575
$ws.Range("N71").Value = -21712.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4999.8887
$ws.Range("I62").Value = 4999.8887
$ws.Range("K62").Value = 4999.8887
$ws.Range("M62").Value = -4375.8887

$ws.Range("H65").Value = 4999.8887
$ws.Range("I65").Value = 4999.8887
$ws.Range("K65").Value = 24999.4435
$ws.Range("M65").Value = -21879.4435

$ws.Range("H113").Value = 575.2381
$ws.Range("I113").Value = 398.85715
$ws.Range("J113").Value = 928
$ws.Range("K113").Value = 1196.57145
$ws.Range("L113").Value = 2784
$ws.Range("M113").Value = 973.4285500000001
$ws.Range("N113").Value = -7124

$ws.Range("H126").Value = 1681.6666
$ws.Range("I126").Value = 1681.6666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5044.9998
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2574.9998
$ws.Range("N126").ClearContents()
